$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A63:C63").Style = $ws.Range("A62:C62").Style

$ws.Range("A63").NumberFormat = "@"
$ws.Range("A63").Value = "2026/01/12"
$ws.Range("B63").Value = "逃离鸭科夫"
$ws.Range("C63").Value = 1143
